$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UCT1")

$ws.Range("G20").Value = 7.48
$ws.Range("H20").Value = 24.9084
$ws.Range("I20").Value = 49.816800000000001
$ws.Range("J20").Value = 76.260000000000005
$ws.Range("K20").Value = 104.16
$ws.Range("L20").Value = 132.06
$ws.Range("M20").Value = 160
$ws.Range("U20").Value = 7.48
$ws.Range("V20").Value = 24.9084
$ws.Range("W20").Value = 49.816800000000001
$ws.Range("X20").Value = 76.260000000000005
$ws.Range("Y20").Value = 104.16
$ws.Range("Z20").Value = 132.06
$ws.Range("AA20").Value = 160

$ws.Range("G21").Value = 7.48
$ws.Range("H21").Value = 24.9084
$ws.Range("I21").Value = 49.816800000000001
$ws.Range("J21").Value = 76.260000000000005
$ws.Range("K21").Value = 104.16
$ws.Range("L21").Value = 132.06
$ws.Range("M21").Value = 160
$ws.Range("U21").Value = 7.48
$ws.Range("V21").Value = 24.9084
$ws.Range("W21").Value = 49.816800000000001
$ws.Range("X21").Value = 76.260000000000005
$ws.Range("Y21").Value = 104.16
$ws.Range("Z21").Value = 132.06
$ws.Range("AA21").Value = 160

$ws.Range("G22").Value = 110.6292
$ws.Range("H22").Value = 368.39523600000001
$ws.Range("I22").Value = 736.79047200000002
$ws.Range("J22").Value = 1127.8854000000001
$ws.Range("K22").Value = 1540.5264
$ws.Range("L22").Value = 1953.1674
$ws.Range("M22").Value = 2366.4
$ws.Range("U22").Value = 110.6292
$ws.Range("V22").Value = 368.39523600000001
$ws.Range("W22").Value = 736.79047200000002
$ws.Range("X22").Value = 1127.8854000000001
$ws.Range("Y22").Value = 1540.5264
$ws.Range("Z22").Value = 1953.1674
$ws.Range("AA22").Value = 2366.4

$ws.Range("G23").Value = 0.374
$ws.Range("H23").Value = 1.2454200000000002
$ws.Range("I23").Value = 2.4908400000000004
$ws.Range("J23").Value = 3.8130000000000002
$ws.Range("K23").Value = 5.2080000000000002
$ws.Range("L23").Value = 6.6029999999999998
$ws.Range("M23").Value = 8
$ws.Range("U23").Value = 0.374
$ws.Range("V23").Value = 1.2454200000000002
$ws.Range("W23").Value = 2.4908400000000004
$ws.Range("X23").Value = 3.8130000000000002
$ws.Range("Y23").Value = 5.2080000000000002
$ws.Range("Z23").Value = 6.6029999999999998
$ws.Range("AA23").Value = 8

$ws.Range("G24").Value = 0.374
$ws.Range("H24").Value = 1.2454200000000002
$ws.Range("I24").Value = 2.4908400000000004
$ws.Range("J24").Value = 3.8130000000000002
$ws.Range("K24").Value = 5.2080000000000002
$ws.Range("L24").Value = 6.6029999999999998
$ws.Range("M24").Value = 8
$ws.Range("U24").Value = 0.374
$ws.Range("V24").Value = 1.2454200000000002
$ws.Range("W24").Value = 2.4908400000000004
$ws.Range("X24").Value = 3.8130000000000002
$ws.Range("Y24").Value = 5.2080000000000002
$ws.Range("Z24").Value = 6.6029999999999998
$ws.Range("AA24").Value = 8

$ws.Range("G25").Value = 0.374
$ws.Range("H25").Value = 1.2454200000000002
$ws.Range("I25").Value = 2.4908400000000004
$ws.Range("J25").Value = 3.8130000000000002
$ws.Range("K25").Value = 5.2080000000000002
$ws.Range("L25").Value = 6.6029999999999998
$ws.Range("M25").Value = 8
$ws.Range("U25").Value = 0.374
$ws.Range("V25").Value = 1.2454200000000002
$ws.Range("W25").Value = 2.4908400000000004
$ws.Range("X25").Value = 3.8130000000000002
$ws.Range("Y25").Value = 5.2080000000000002
$ws.Range("Z25").Value = 6.6029999999999998
$ws.Range("AA25").Value = 8

$ws.Range("G26").Value = 0.374
$ws.Range("H26").Value = 1.2454200000000002
$ws.Range("I26").Value = 2.4908400000000004
$ws.Range("J26").Value = 3.8130000000000002
$ws.Range("K26").Value = 5.2080000000000002
$ws.Range("L26").Value = 6.6029999999999998
$ws.Range("M26").Value = 8
$ws.Range("U26").Value = 0.374
$ws.Range("V26").Value = 1.2454200000000002
$ws.Range("W26").Value = 2.4908400000000004
$ws.Range("X26").Value = 3.8130000000000002
$ws.Range("Y26").Value = 5.2080000000000002
$ws.Range("Z26").Value = 6.6029999999999998
$ws.Range("AA26").Value = 8

$ws.Range("G27").Value = 0.374
$ws.Range("H27").Value = 1.2454200000000002
$ws.Range("I27").Value = 2.4908400000000004
$ws.Range("J27").Value = 3.8130000000000002
$ws.Range("K27").Value = 5.2080000000000002
$ws.Range("L27").Value = 6.6029999999999998
$ws.Range("M27").Value = 8
$ws.Range("U27").Value = 0.374
$ws.Range("V27").Value = 1.2454200000000002
$ws.Range("W27").Value = 2.4908400000000004
$ws.Range("X27").Value = 3.8130000000000002
$ws.Range("Y27").Value = 5.2080000000000002
$ws.Range("Z27").Value = 6.6029999999999998
$ws.Range("AA27").Value = 8

$ws.Range("G28").Value = 0.374
$ws.Range("H28").Value = 1.2454200000000002
$ws.Range("I28").Value = 2.4908400000000004
$ws.Range("J28").Value = 3.8130000000000002
$ws.Range("K28").Value = 5.2080000000000002
$ws.Range("L28").Value = 6.6029999999999998
$ws.Range("M28").Value = 8
$ws.Range("U28").Value = 0.374
$ws.Range("V28").Value = 1.2454200000000002
$ws.Range("W28").Value = 2.4908400000000004
$ws.Range("X28").Value = 3.8130000000000002
$ws.Range("Y28").Value = 5.2080000000000002
$ws.Range("Z28").Value = 6.6029999999999998
$ws.Range("AA28").Value = 8

$ws.Range("G29").Value = 0.374
$ws.Range("H29").Value = 1.2454200000000002
$ws.Range("I29").Value = 2.4908400000000004
$ws.Range("J29").Value = 3.8130000000000002
$ws.Range("K29").Value = 5.2080000000000002
$ws.Range("L29").Value = 6.6029999999999998
$ws.Range("M29").Value = 8
$ws.Range("U29").Value = 0.374
$ws.Range("V29").Value = 1.2454200000000002
$ws.Range("W29").Value = 2.4908400000000004
$ws.Range("X29").Value = 3.8130000000000002
$ws.Range("Y29").Value = 5.2080000000000002
$ws.Range("Z29").Value = 6.6029999999999998
$ws.Range("AA29").Value = 8

$ws.Range("G30").Value = 0.374
$ws.Range("H30").Value = 1.2454200000000002
$ws.Range("I30").Value = 2.4908400000000004
$ws.Range("J30").Value = 3.8130000000000002
$ws.Range("K30").Value = 5.2080000000000002
$ws.Range("L30").Value = 6.6029999999999998
$ws.Range("M30").Value = 8
$ws.Range("U30").Value = 0.374
$ws.Range("V30").Value = 1.2454200000000002
$ws.Range("W30").Value = 2.4908400000000004
$ws.Range("X30").Value = 3.8130000000000002
$ws.Range("Y30").Value = 5.2080000000000002
$ws.Range("Z30").Value = 6.6029999999999998
$ws.Range("AA30").Value = 8

$ws.Range("G31").Value = 0.374
$ws.Range("H31").Value = 1.2454200000000002
$ws.Range("I31").Value = 2.4908400000000004
$ws.Range("J31").Value = 3.8130000000000002
$ws.Range("K31").Value = 5.2080000000000002
$ws.Range("L31").Value = 6.6029999999999998
$ws.Range("M31").Value = 8
$ws.Range("U31").Value = 0.374
$ws.Range("V31").Value = 1.2454200000000002
$ws.Range("W31").Value = 2.4908400000000004
$ws.Range("X31").Value = 3.8130000000000002
$ws.Range("Y31").Value = 5.2080000000000002
$ws.Range("Z31").Value = 6.6029999999999998
$ws.Range("AA31").Value = 8

$ws.Range("G32").Value = 0.374
$ws.Range("H32").Value = 1.2454200000000002
$ws.Range("I32").Value = 2.4908400000000004
$ws.Range("J32").Value = 3.8130000000000002
$ws.Range("K32").Value = 5.2080000000000002
$ws.Range("L32").Value = 6.6029999999999998
$ws.Range("M32").Value = 8
$ws.Range("U32").Value = 0.374
$ws.Range("V32").Value = 1.2454200000000002
$ws.Range("W32").Value = 2.4908400000000004
$ws.Range("X32").Value = 3.8130000000000002
$ws.Range("Y32").Value = 5.2080000000000002
$ws.Range("Z32").Value = 6.6029999999999998
$ws.Range("AA32").Value = 8

$ws.Range("G33").Value = 0.374
$ws.Range("H33").Value = 1.2454200000000002
$ws.Range("I33").Value = 2.4908400000000004
$ws.Range("J33").Value = 3.8130000000000002
$ws.Range("K33").Value = 5.2080000000000002
$ws.Range("L33").Value = 6.6029999999999998
$ws.Range("M33").Value = 8
$ws.Range("U33").Value = 0.374
$ws.Range("V33").Value = 1.2454200000000002
$ws.Range("W33").Value = 2.4908400000000004
$ws.Range("X33").Value = 3.8130000000000002
$ws.Range("Y33").Value = 5.2080000000000002
$ws.Range("Z33").Value = 6.6029999999999998
$ws.Range("AA33").Value = 8

$ws.Range("G34").Value = 0.374
$ws.Range("H34").Value = 1.2454200000000002
$ws.Range("I34").Value = 2.4908400000000004
$ws.Range("J34").Value = 3.8130000000000002
$ws.Range("K34").Value = 5.2080000000000002
$ws.Range("L34").Value = 6.6029999999999998
$ws.Range("M34").Value = 8
$ws.Range("U34").Value = 0.374
$ws.Range("V34").Value = 1.2454200000000002
$ws.Range("W34").Value = 2.4908400000000004
$ws.Range("X34").Value = 3.8130000000000002
$ws.Range("Y34").Value = 5.2080000000000002
$ws.Range("Z34").Value = 6.6029999999999998
$ws.Range("AA34").Value = 8

$ws.Range("G35").Value = 21.692
$ws.Range("H35").Value = 72.234359999999995
$ws.Range("I35").Value = 144.46871999999999
$ws.Range("J35").Value = 221.154
$ws.Range("K35").Value = 302.06399999999996
$ws.Range("L35").Value = 382.97399999999999
$ws.Range("M35").Value = 463.99999999999994
$ws.Range("U35").Value = 21.692
$ws.Range("V35").Value = 72.234359999999995
$ws.Range("W35").Value = 144.46871999999999
$ws.Range("X35").Value = 221.154
$ws.Range("Y35").Value = 302.06399999999996
$ws.Range("Z35").Value = 382.97399999999999
$ws.Range("AA35").Value = 463.99999999999994

$ws.Range("G36").Value = 9.35
$ws.Range("H36").Value = 31.135500000000004
$ws.Range("I36").Value = 62.271000000000008
$ws.Range("J36").Value = 95.325000000000003
$ws.Range("K36").Value = 130.20000000000002
$ws.Range("L36").Value = 165.07500000000002
$ws.Range("M36").Value = 200
$ws.Range("U36").Value = 9.35
$ws.Range("V36").Value = 31.135500000000004
$ws.Range("W36").Value = 62.271000000000008
$ws.Range("X36").Value = 95.325000000000003
$ws.Range("Y36").Value = 130.20000000000002
$ws.Range("Z36").Value = 165.07500000000002
$ws.Range("AA36").Value = 200

$ws.Range("G37").Value = 0.074800000000000005373479439186
$ws.Range("H37").Value = 0.24908400000000003
$ws.Range("I37").Value = 0.49816800000000006
$ws.Range("J37").Value = 0.76260000000000006
$ws.Range("K37").Value = 1.0416000000000001
$ws.Range("L37").Value = 1.3206
$ws.Range("M37").Value = 1.6
$ws.Range("U37").Value = 0.074800000000000005373479439186
$ws.Range("V37").Value = 0.24908400000000003
$ws.Range("W37").Value = 0.49816800000000006
$ws.Range("X37").Value = 0.76260000000000006
$ws.Range("Y37").Value = 1.0416000000000001
$ws.Range("Z37").Value = 1.3206
$ws.Range("AA37").Value = 1.6

$ws.Range("G38").Value = 19.073999999999998
$ws.Range("H38").Value = 63.516419999999997
$ws.Range("I38").Value = 127.03283999999999
$ws.Range("J38").Value = 194.46299999999999
$ws.Range("K38").Value = 265.608
$ws.Range("L38").Value = 336.75299999999999
$ws.Range("M38").Value = 408
$ws.Range("U38").Value = 19.073999999999998
$ws.Range("V38").Value = 63.516419999999997
$ws.Range("W38").Value = 127.03283999999999
$ws.Range("X38").Value = 194.46299999999999
$ws.Range("Y38").Value = 265.608
$ws.Range("Z38").Value = 336.75299999999999
$ws.Range("AA38").Value = 408

$ws.Range("G39").Value = 1.496
$ws.Range("H39").Value = 4.9816800000000008
$ws.Range("I39").Value = 9.9633600000000015
$ws.Range("J39").Value = 15.252000000000001
$ws.Range("K39").Value = 20.832000000000001
$ws.Range("L39").Value = 26.411999999999999
$ws.Range("M39").Value = 32
$ws.Range("U39").Value = 1.496
$ws.Range("V39").Value = 4.9816800000000008
$ws.Range("W39").Value = 9.9633600000000015
$ws.Range("X39").Value = 15.252000000000001
$ws.Range("Y39").Value = 20.832000000000001
$ws.Range("Z39").Value = 26.411999999999999
$ws.Range("AA39").Value = 32

$ws.Range("G40").Value = 2.2440000000000002
$ws.Range("H40").Value = 7.4725200000000003
$ws.Range("I40").Value = 14.945040000000001
$ws.Range("J40").Value = 22.878
$ws.Range("K40").Value = 31.248000000000001
$ws.Range("L40").Value = 39.618000000000002
$ws.Range("M40").Value = 48
$ws.Range("U40").Value = 2.2440000000000002
$ws.Range("V40").Value = 7.4725200000000003
$ws.Range("W40").Value = 14.945040000000001
$ws.Range("X40").Value = 22.878
$ws.Range("Y40").Value = 31.248000000000001
$ws.Range("Z40").Value = 39.618000000000002
$ws.Range("AA40").Value = 48

$ws.Range("G41").Value = 187
$ws.Range("H41").Value = 622.71
$ws.Range("I41").Value = 1245.42
$ws.Range("J41").Value = 1906.5
$ws.Range("K41").Value = 2604
$ws.Range("L41").Value = 3301.5
$ws.Range("M41").Value = 4000
$ws.Range("U41").Value = 187
$ws.Range("V41").Value = 622.71
$ws.Range("W41").Value = 1245.42
$ws.Range("X41").Value = 1906.5
$ws.Range("Y41").Value = 2604
$ws.Range("Z41").Value = 3301.5
$ws.Range("AA41").Value = 4000

$ws.Range("G42").Value = 0.748
$ws.Range("H42").Value = 2.4908400000000004
$ws.Range("I42").Value = 4.9816800000000008
$ws.Range("J42").Value = 7.6260000000000003
$ws.Range("K42").Value = 10.416
$ws.Range("L42").Value = 13.206
$ws.Range("M42").Value = 16
$ws.Range("U42").Value = 0.748
$ws.Range("V42").Value = 2.4908400000000004
$ws.Range("W42").Value = 4.9816800000000008
$ws.Range("X42").Value = 7.6260000000000003
$ws.Range("Y42").Value = 10.416
$ws.Range("Z42").Value = 13.206
$ws.Range("AA42").Value = 16

$ws.Range("G43").Value = 2.2440000000000002
$ws.Range("H43").Value = 7.4725200000000003
$ws.Range("I43").Value = 14.945040000000001
$ws.Range("J43").Value = 22.878
$ws.Range("K43").Value = 31.248000000000001
$ws.Range("L43").Value = 39.618000000000002
$ws.Range("M43").Value = 48
$ws.Range("U43").Value = 2.2440000000000002
$ws.Range("V43").Value = 7.4725200000000003
$ws.Range("W43").Value = 14.945040000000001
$ws.Range("X43").Value = 22.878
$ws.Range("Y43").Value = 31.248000000000001
$ws.Range("Z43").Value = 39.618000000000002
$ws.Range("AA43").Value = 48

$ws.Range("G20:M43").Select()